# Auto update Excel log
# Appends newly-logged sensor readings to the PIR, Humidity, Temperature and
# mmWave sheets (matching the device export that ran on 2026-01-28).

$wb = $excel.ActiveWorkbook

function Append-LogRow {
    param($ws, [int]$Row, [string]$Date, [string]$Time, [string]$Hour, [string]$Location, [string]$Value, [string]$Status, [bool]$ValueIsPercentLike = $false)

    # Column A holds a date-shaped string ("2026-01-28"); briefly force text
    # so it is stored the same way as every other row in the log instead of
    # being auto-converted into a real date serial number, then restore the
    # General format so the cell keeps the sheet's default appearance.
    $cellA = $ws.Cells.Item($Row, 1)
    $cellA.NumberFormat = "@"
    $cellA.Value = $Date
    $cellA.NumberFormat = "General"

    $ws.Cells.Item($Row, 2).Value = $Time
    $ws.Cells.Item($Row, 3).Value = $Hour
    $ws.Cells.Item($Row, 4).Value = $Location

    # The "Value" column is free-form text everywhere in the log, but a
    # percentage-shaped reading ("88.2%") would otherwise get silently
    # converted into a real numeric percentage - force text only for that
    # case so everything else keeps its natural (already-text) type, then
    # restore General the same way as above.
    $cellE = $ws.Cells.Item($Row, 5)
    if ($ValueIsPercentLike) {
        $cellE.NumberFormat = "@"
        $cellE.Value = $Value
        $cellE.NumberFormat = "General"
    } else {
        $cellE.Value = $Value
    }

    $ws.Cells.Item($Row, 6).Value = $Status
}

# ---------------------------------------------------------------------------
# PIR sheet - rows 161-174 (Bathroom / No Motion / Inactive)
# ---------------------------------------------------------------------------
$wsPIR = $wb.Worksheets.Item("PIR")
$pirTimes = @(
    "15:01:01","15:01:02","15:01:05","15:01:10","15:01:16","15:01:20",
    "15:01:25","15:01:30","15:01:36","15:01:40","15:01:45","15:01:50",
    "15:01:56","15:02:00"
)
$row = 161
foreach ($t in $pirTimes) {
    Append-LogRow $wsPIR $row "2026-01-28" $t "15:00" "Bathroom" "No Motion" "Inactive"
    $row++
}

# ---------------------------------------------------------------------------
# Humidity sheet - rows 155-168 (Bathroom / percentage / Active)
# ---------------------------------------------------------------------------
$wsHumidity = $wb.Worksheets.Item("Humidity")
$humidityData = @(
    @("15:01:01","88.2%"),
    @("15:01:03","88.3%"),
    @("15:01:06","87.3%"),
    @("15:01:14","88.3%"),
    @("15:01:18","87.4%"),
    @("15:01:22","88.3%"),
    @("15:01:31","88.3%"),
    @("15:01:35","88.3%"),
    @("15:01:39","87.4%"),
    @("15:01:43","88.3%"),
    @("15:01:47","87.4%"),
    @("15:01:51","88.3%"),
    @("15:01:55","88.3%"),
    @("15:01:59","87.4%")
)
$row = 155
foreach ($item in $humidityData) {
    Append-LogRow $wsHumidity $row "2026-01-28" $item[0] "15:00" "Bathroom" $item[1] "Active" $true
    $row++
}

# ---------------------------------------------------------------------------
# Temperature sheet - rows 155-168 (Bathroom / Celsius / Active)
# ---------------------------------------------------------------------------
$wsTemperature = $wb.Worksheets.Item("Temperature")
$temperatureData = @(
    @("15:01:02","22.9C"),
    @("15:01:03","22.9C"),
    @("15:01:07","22.8C"),
    @("15:01:15","22.9C"),
    @("15:01:19","22.9C"),
    @("15:01:23","22.9C"),
    @("15:01:31","22.9C"),
    @("15:01:35","22.9C"),
    @("15:01:39","22.8C"),
    @("15:01:43","22.9C"),
    @("15:01:47","22.9C"),
    @("15:01:51","22.9C"),
    @("15:01:55","22.9C"),
    @("15:01:59","22.9C")
)
$row = 155
foreach ($item in $temperatureData) {
    Append-LogRow $wsTemperature $row "2026-01-28" $item[0] "15:00" "Bathroom" $item[1] "Active"
    $row++
}

# ---------------------------------------------------------------------------
# mmWave sheet - row 6 (Living Room / No Presence / Inactive)
# ---------------------------------------------------------------------------
$wsMmWave = $wb.Worksheets.Item("mmWave")
Append-LogRow $wsMmWave 6 "2026-01-28" "15:01:44" "15:00" "Living Room" "No Presence" "Inactive"
